# Actualización desde MV -datos-
# Append the new daily bond-rate observations (16..29 Sep 2021) below the
# existing table, which currently ends at row 177.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, date (col A), col C, col D, col E ($null = leave blank)
$rows = @(
    @(178, "16-09-2021", 3.8,   $null, 5.17),
    @(179, "20-09-2021", 3.85,  4.73,  5.15),
    @(180, "21-09-2021", 3.8,   $null, 5.01),
    @(181, "22-09-2021", $null, $null, 5.02),
    @(182, "23-09-2021", 4,     4.7,   5.15),
    @(183, "24-09-2021", 4.07,  4.85,  5.24),
    @(184, "27-09-2021", 4.1,   5,     5.4),
    @(185, "28-09-2021", 4.1,   5.1,   5.5),
    @(186, "29-09-2021", 4.4,   5.18,  5.56)
)

foreach ($entry in $rows) {
    $r    = $entry[0]
    $date = $entry[1]
    $c    = $entry[2]
    $d    = $entry[3]
    $e    = $entry[4]

    $ws.Cells.Item($r, 1).Value = $date
    if ($c -ne $null) { $ws.Cells.Item($r, 3).Value = $c }
    if ($d -ne $null) { $ws.Cells.Item($r, 4).Value = $d }
    if ($e -ne $null) { $ws.Cells.Item($r, 5).Value = $e }
}
